# Apply "Natmi following Dr Hou advice" update:
# Sending-cluster "ECs" is added as a third cluster (alongside existing FAPs and sCs)
# for the Lgi2 -> Adam22 ligand-receptor pair, expanding the 2x3 grid of
# (sending cluster x target cluster) combinations into a full 3x3 grid (rows 2-10),
# and all NATMI-derived statistic columns (E:T) are recalculated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Lgi2", "Adam22", "ECs", 1, 0.3333333333333333, 0.058936, 0.176808, 0.003640699631737656, 0.003640699631737656, 3, 1, 0.8328803333333332, 2.498641, 0.03917234793046614, 0.03917234793046613, 0.04908663532533333, 0.4417797179279999, 0.0001426147526847474, 0.0001426147526847474),
    @("ECs", "Lgi2", "Adam22", "FAPs", 1, 0.3333333333333333, 0.058936, 0.176808, 0.003640699631737656, 0.003640699631737656, 3, 1, 4.873409, 14.620227, 0.2292080450398417, 0.2292080450398417, 0.287219232824, 2.584973095416, 0.0008344776451678599, 0.0008344776451678599),
    @("ECs", "Lgi2", "Adam22", "sCs", 1, 0.3333333333333333, 0.058936, 0.176808, 0.003640699631737656, 0.003640699631737656, 3, 1, 15.55565633333333, 46.66696899999999, 0.7316196070296923, 0.7316196070296922, 0.9167881616613331, 8.251093454951999, 0.002663607233885049, 0.002663607233885049),
    @("FAPs", "Lgi2", "Adam22", "ECs", 3, 1, 14.32797533333333, 42.983926, 0.8850932285803735, 0.8850932285803734, 3, 1, 0.8328803333333332, 2.498641, 0.03917234793046614, 0.03917234793046613, 11.93348887161844, 107.401399844566, 0.03467117990084999, 0.03467117990084997),
    @("FAPs", "Lgi2", "Adam22", "FAPs", 3, 1, 14.32797533333333, 42.983926, 0.8850932285803735, 0.8850932285803734, 3, 1, 4.873409, 14.620227, 0.2292080450398417, 0.2292080450398417, 69.82608394124466, 628.4347554712019, 0.2028704886009092, 0.2028704886009091),
    @("FAPs", "Lgi2", "Adam22", "sCs", 3, 1, 14.32797533333333, 42.983926, 0.8850932285803735, 0.8850932285803734, 3, 1, 15.55565633333333, 46.66696899999999, 0.7316196070296923, 0.7316196070296922, 222.8810602378104, 2005.929542140294, 0.6475515600786145, 0.6475515600786144),
    @("sCs", "Lgi2", "Adam22", "ECs", 3, 1, 1.801186, 5.403558, 0.111266071787889, 0.1112660717878889, 3, 1, 0.8328803333333332, 2.498641, 0.03917234793046614, 0.03917234793046613, 1.500172396075333, 13.501551564678, 0.004358553276931409, 0.004358553276931408),
    @("sCs", "Lgi2", "Adam22", "FAPs", 3, 1, 1.801186, 5.403558, 0.111266071787889, 0.1112660717878889, 3, 1, 4.873409, 14.620227, 0.2292080450398417, 0.2292080450398417, 8.777916063074, 79.001244567666, 0.02550307879376471, 0.02550307879376471),
    @("sCs", "Lgi2", "Adam22", "sCs", 3, 1, 1.801186, 5.403558, 0.111266071787889, 0.1112660717878889, 3, 1, 15.55565633333333, 46.66696899999999, 0.7316196070296923, 0.7316196070296922, 28.01863040841133, 252.167673675702, 0.08140443971719286, 0.08140443971719283)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $row[$j]
    }
}